# Mini Project 1 Completed
#
# The "Signup" worksheet has a small table of signup test-data in A1:F5.
# Row 2 held a test case whose first name was "Hate" and whose e-mail
# address was "hatebrotest@gmail.com". This updates that row to use
# "Hope" / "HopenNopeh@gmail.com" instead, and also gives the e-mail
# address in C3 the same (hyperlink) formatting already used by the
# other e-mail cells in the column (C2 and C5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Signup")

# Row 2: firstname "Hate" -> "Hope", email "hatebrotest@gmail.com" -> "HopenNopeh@gmail.com"
$ws.Range("A2").Value = "Hope"
$ws.Range("C2").Value = "HopenNopeh@gmail.com"

# Give C3 (another e-mail/hyperlink cell) the same formatting as C2,
# matching the hyperlink-style look already applied to C2 and C5.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
